$wb = $excel.ActiveWorkbook
$wsExhibition = $wb.Worksheets.Item("展览")
$wsAllTypes = $wb.Worksheets.Item("全部类型")

# Update "想去人数" (F column) values on 展览 sheet
$wsExhibition.Range("F2").Value = 79
$wsExhibition.Range("F3").Value = 411
$wsExhibition.Range("F4").Value = 3012
$wsExhibition.Range("F6").Value = 36
$wsExhibition.Range("F9").Value = 4
$wsExhibition.Range("F10").Value = 14353
$wsExhibition.Range("F11").Value = 150
$wsExhibition.Range("F12").Value = 117
$wsExhibition.Range("F13").Value = 5739
$wsExhibition.Range("F15").Value = 68
$wsExhibition.Range("F16").Value = 41
$wsExhibition.Range("F17").Value = 58
$wsExhibition.Range("F18").Value = 1233
$wsExhibition.Range("F20").Value = 78
$wsExhibition.Range("F22").Value = 786
$wsExhibition.Range("F23").Value = 2928
$wsExhibition.Range("F25").Value = 10531
$wsExhibition.Range("F26").Value = 1199
$wsExhibition.Range("F27").Value = 55
$wsExhibition.Range("F29").Value = 3732
$wsExhibition.Range("F31").Value = 65

# Update "想去人数" (F column) values on 全部类型 sheet
$wsAllTypes.Range("F2").Value = 79
$wsAllTypes.Range("F3").Value = 411
$wsAllTypes.Range("F5").Value = 3012
$wsAllTypes.Range("F7").Value = 36
$wsAllTypes.Range("F10").Value = 4
$wsAllTypes.Range("F11").Value = 14353
$wsAllTypes.Range("F12").Value = 150
$wsAllTypes.Range("F13").Value = 117
$wsAllTypes.Range("F14").Value = 5739
$wsAllTypes.Range("F16").Value = 68
$wsAllTypes.Range("F17").Value = 41
$wsAllTypes.Range("F18").Value = 58
$wsAllTypes.Range("F19").Value = 1233
$wsAllTypes.Range("F21").Value = 78
$wsAllTypes.Range("F23").Value = 786
$wsAllTypes.Range("F24").Value = 2928
$wsAllTypes.Range("F27").Value = 10531
$wsAllTypes.Range("F28").Value = 1199
$wsAllTypes.Range("F29").Value = 55
$wsAllTypes.Range("F31").Value = 3732
$wsAllTypes.Range("F33").Value = 65

